$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.695.02"
$ws.Range("E2").Value = "'  -0.23%  "
$ws.Range("D3").Value = "'3.581.24"
$ws.Range("E3").Value = "'  -2.29%  "
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'237.63"
$ws.Range("E5").Value = "'  -1.89%  "
$ws.Range("D6").Value = "'654.10"
$ws.Range("E6").Value = "'  +1.36%  "
$ws.Range("D7").Value = "'1.52"
$ws.Range("E7").Value = "'  +3.72%  "
$ws.Range("D8").Value = "'0.402"
$ws.Range("E8").Value = "'  +0.13%  "
$ws.Range("E9").Value = "'  +0.08%  "
$ws.Range("E10").Value = "'  +1.73%  "
$ws.Range("D11").Value = "'3.581.77"
$ws.Range("E11").Value = "'  -2.25%  "
$ws.Range("E12").Value = "'  +1.11%  "
$ws.Range("D13").Value = "'42.87"
$ws.Range("E13").Value = "'  -2.33%  "
$ws.Range("D14").Value = "'6.45"
$ws.Range("E14").Value = "'  +0.98%  "
$ws.Range("D15").Value = "'4.249.74"
$ws.Range("E15").Value = "'  -2.60%  "
$ws.Range("D16").Value = "'95.542.43"
$ws.Range("E16").Value = "'  -0.29%  "
$ws.Range("E17").Value = "'  -0.89%  "
$ws.Range("D18").Value = "'3.571.40"
$ws.Range("E18").Value = "'  -2.64%  "
$ws.Range("D19").Value = "'12.75"
$ws.Range("E19").Value = "'  -5.00%  "
$ws.Range("D20").Value = "'7.75"
$ws.Range("E20").Value = "'  -3.84%  "
$ws.Range("D21").Value = "'17.96"
$ws.Range("E21").Value = "'  -3.49%  "
$ws.Range("D22").Value = "'0.493"
$ws.Range("E22").Value = "'  +2.64%  "
$ws.Range("D23").Value = "'3.44"
$ws.Range("E23").Value = "'  +0.37%  "
$ws.Range("D24").Value = "'511.03"
$ws.Range("E24").Value = "'  -1.67%  "
$ws.Range("D25").Value = "'7.08"
$ws.Range("E25").Value = "'  +3.55%  "
$ws.Range("D26").Value = "'0.0000198"
$ws.Range("E26").Value = "'  +0.31%  "
$ws.Range("D27").Value = "'95.89"
$ws.Range("E27").Value = "'  +2.32%  "
$ws.Range("D28").Value = "'12.82"
$ws.Range("E28").Value = "'  +1.16%  "
$ws.Range("D29").Value = "'3.774.35"
$ws.Range("E29").Value = "'  -2.43%  "
$ws.Range("D30").Value = "'3.04"
$ws.Range("E30").Value = "'  -3.41%  "
$ws.Range("E31").Value = "'  +2.70%  "
$ws.Range("D32").Value = "'11.57"
$ws.Range("E32").Value = "'  -1.12%  "
$ws.Range("E33").Value = "'  +0.15%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "'  +0.08%  "
$ws.Range("E35").Value = "'  -1.10%  "
$ws.Range("D36").Value = "'31.96"
$ws.Range("E36").Value = "'  -3.18%  "
$ws.Range("D37").Value = "'1.69"
$ws.Range("E37").Value = "'  +12.18%  "
$ws.Range("D38").Value = "'0.564"
$ws.Range("E38").Value = "'  -2.58%  "
$ws.Range("D39").Value = "'8.60"
$ws.Range("E39").Value = "'  +8.29%  "
$ws.Range("D40").Value = "'596.16"
$ws.Range("E40").Value = "'  +6.13%  "
$ws.Range("E41").Value = "'  -1.01%  "
$ws.Range("E42").Value = "'  +0.09%  "
$ws.Range("D43").Value = "'1.87"
$ws.Range("E43").Value = "'  +6.79%  "
$ws.Range("D44").Value = "'0.913"
$ws.Range("E44").Value = "'  -5.60%  "
$ws.Range("E45").Value = "'  +0.06%  "
$ws.Range("E46").Value = "'  +3.90%  "
$ws.Range("D47").Value = "'34.32"
$ws.Range("E47").Value = "'  +1.54%  "
$ws.Range("D48").Value = "'23.45"
$ws.Range("E48").Value = "'  -1.07%  "
$ws.Range("E49").Value = "'  -2.83%  "
$ws.Range("D50").Value = "'3.49"
$ws.Range("E50").Value = "'  -0.05%  "
$ws.Range("D51").Value = "'8.25"
$ws.Range("E51").Value = "'  -0.69%  "
